$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header labels: drop the leading "%" sign from the IVA / ICE columns
$ws.Range("M10").Value = "IVA"
$ws.Range("N10").Value = "ICE"

# Move the active selection as it was left by the author
$ws.Range("N11").Select()
